# Fix figure displacement: the "GalaxyUserCount.png" picture together with
# its six annotation shapes (two "stem" connectors, two arrow connectors and
# the two "First/Second Galaxy Workshop" labels) are combined into a single
# group, and the annotation shapes are repositioned so they line up with the
# picture again.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the picture and the six annotation shapes by their (stable) shape Id,
# regardless of their current index in the shape collection.
function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

$picShape   = Get-ShapeById $s.Shapes 19   # Bild 18 - GalaxyUserCount.png
$cxn23      = Get-ShapeById $s.Shapes 23   # Gerade Verbindung 22
$cxn45      = Get-ShapeById $s.Shapes 45   # Gerade Verbindung 44
$cxn51      = Get-ShapeById $s.Shapes 51   # Gerade Verbindung mit Pfeil 50
$sp46       = Get-ShapeById $s.Shapes 46   # Textfeld 45 (Second Galaxy Workshop)
$cxn56      = Get-ShapeById $s.Shapes 56   # Gerade Verbindung mit Pfeil 55
$sp57       = Get-ShapeById $s.Shapes 57   # Textfeld 56 (First Galaxy Workshop)

# Build the range in the desired final (z-)order and group the shapes.
$idxs = @($picShape.Id, $cxn23.Id, $cxn45.Id, $cxn51.Id, $sp46.Id, $cxn56.Id, $sp57.Id) | ForEach-Object {
    $shp = Get-ShapeById $s.Shapes $_
    $shp.ZOrderPosition
}
$range = $s.Shapes.Range($idxs)
$grp = $range.Group()

# Put the new group where the picture used to live (the very back of the
# z-order / first shape in the XML) and give it its proper (localized) name.
$grp.ZOrder(1)
$grp.Name = "Gruppierung 2"

# Re-align the annotation shapes with the picture (new absolute positions,
# sizes stay the same). Values are EMU/12700 = points, as used by the Shape
# Left/Top COM properties.
$positions = @{
    23 = @(1904.5763779527558, 2813.7711811023623)  # Gerade Verbindung 22
    45 = @(2002.499842519685,  2627.2667716535434)  # Gerade Verbindung 44
    51 = @(1975.999842519685,  2673.0377952755907)  # Gerade Verbindung mit Pfeil 50
    46 = @(1836.999842519685,  2625.780787401575)   # Textfeld 45
    56 = @(1876.0,             2861.028188976378)   # Gerade Verbindung mit Pfeil 55
    57 = @(1737.0,             2813.7711811023623)  # Textfeld 56
}

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $item = $grp.GroupItems.Item($i)
    if ($positions.ContainsKey($item.Id)) {
        $xy = $positions[$item.Id]
        $item.Left = $xy[0]
        $item.Top = $xy[1]
    }
}
